$d = $word.ActiveDocument

$replacements = @(
    @("54-34=", "36-18="),
    @("26+65=", "2+78="),
    @("72+21=", "38-1="),
    @("45+50=", "38-22="),
    @("43+25=", "17+29="),
    @("99-18=", "41+52="),
    @("66+31=", "58+10="),
    @("7+1=", "18+30="),
    @("89-84=", "33+61="),
    @("66-26=", "7+69="),
    @("10-7=", "26+62="),
    @("30+59=", "94-46="),
    @("55-8=", "32+44="),
    @("38-25=", "17+17="),
    @("34+30=", "99-78="),
    @("59-43=", "28+14="),
    @("18+70=", "86+5="),
    @("79-48=", "46+37="),
    @("93-27=", "73-37="),
    @("69-32=", "94-76="),
    @("10-0=", "30+36="),
    @("47+4=", "57+8="),
    @("14+24=", "30+65="),
    @("58+1=", "36+12="),
    @("61-52=", "62-20="),
    @("50+14=", "41-15="),
    @("7+22=", "10+72="),
    @("43-19=", "52-9="),
    @("66+13=", "56-26="),
    @("54+15=", "81-48="),
    @("24-9=", "88+6="),
    @("9+60=", "58-40="),
    @("49+42=", "68-17="),
    @("30+26=", "6+62="),
    @("93-52=", "85-4="),
    @("20+50=", "63-45="),
    @("8+4=", "5+80="),
    @("45-9=", "88-44="),
    @("36+57=", "81-16="),
    @("76-26=", "81-24="),
    @("21+4=", "76-32="),
    @("79-30=", "57-19="),
    @("61+1=", "86+5="),
    @("13-1=", "53+29="),
    @("72-49=", "91-25="),
    @("97-0=", "55+28="),
    @("69-18=", "85-30="),
    @("40-22=", "93-91="),
    @("73-31=", "92-63="),
    @("27+66=", "99-81="),
    @("31-22=", "40-35="),
    @("68-26=", "51+18="),
    @("30+4=", "77-76="),
    @("54+6=", "89-21="),
    @("19+16=", "5+52="),
    @("32+6=", "71-44="),
    @("3+57=", "8+68="),
    @("55+26=", "70-29="),
    @("4+45=", "65+8="),
    @("96-82=", "36+28="),
    @("11+74=", "26+60="),
    @("52-29=", "55+12="),
    @("7+52=", "82-59="),
    @("47+13=", "15+62="),
    @("45-2=", "34+50="),
    @("82-81=", "41+24="),
    @("88-83=", "37-24="),
    @("33+11=", "89-7="),
    @("35+29=", "42+50="),
    @("55+38=", "90-75="),
    @("17+61=", "24+13="),
    @("77+20=", "88-51="),
    @("96-95=", "57+33="),
    @("84+2=", "31-11="),
    @("53+38=", "62+16="),
    @("28-8=", "95-44="),
    @("71-65=", "73+11="),
    @("25+18=", "87-80="),
    @("73+3=", "9+84="),
    @("9+67=", "59+40="),
    @("5+3=", "97-24="),
    @("30+49=", "51-40="),
    @("70-10=", "77-63="),
    @("87+7=", "23-3="),
    @("36-26=", "5+66="),
    @("47-2=", "85-28="),
    @("80-44=", "96-61="),
    @("43+7=", "92+1="),
    @("14+10=", "33-8="),
    @("15+56=", "50-2="),
    @("27+44=", "35+35="),
    @("83-1=", "56-9="),
    @("64-0=", "14-9="),
    @("82-58=", "96+2="),
    @("65-27=", "64+17="),
    @("81-51=", "1+25="),
    @("79-60=", "40+14="),
    @("76-2=", "54+1="),
    @("71+22=", "92-84="),
    @("57-38=", "5+14=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
